$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Department Name" header becomes "Department"
$ws.Range("H1").Value = "Department"

# New Department Code / Department columns (G/H) populated for each data row
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 7).Value = "CMU"
    $ws.Cells.Item($row, 8).Value = "Khoa đào tạo quốc tế"
}
